# Auto-generated edit script applying the row-permutation diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$ws.Range("A2").Value = 130886842
$ws.Range("Q2").Value = 434316
$ws.Range("R2").Value = 7052462
$ws.Range("A3").Value = 130886843
$ws.Range("Q3").Value = 434321
$ws.Range("R3").Value = 7052458
$ws.Range("A11").Value = 130886779
$ws.Range("B11").Value = 57884
$ws.Range("E11").Value = 100109
$ws.Range("F11").Value = "Tretåig hackspett"
$ws.Range("G11").Value = "Picoides tridactylus"
$ws.Range("H11").Value = "(Linnaeus, 1758)"
$ws.Range("Q11").Value = 434359
$ws.Range("R11").Value = 7052173
$ws.Range("AC11").Value = "Ringhack äldre"
$ws.Range("A12").Value = 130886838
$ws.Range("B12").Value = 91828
$ws.Range("E12").Value = 5432
$ws.Range("F12").Value = "Granticka"
$ws.Range("G12").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H12").Value = ""
$ws.Range("Q12").Value = 434362
$ws.Range("R12").Value = 7052349
$ws.Range("AC12").Value = $null
$ws.Range("A21").Value = 130886836
$ws.Range("B21").Value = 91828
$ws.Range("E21").Value = 5432
$ws.Range("F21").Value = "Granticka"
$ws.Range("G21").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H21").Value = ""
$ws.Range("Q21").Value = 434854
$ws.Range("R21").Value = 7051718
$ws.Range("AC21").Value = $null
$ws.Range("A22").Value = 130886792
$ws.Range("Q22").Value = 434158
$ws.Range("R22").Value = 7052168
$ws.Range("AC22").Value = "Ringhack äldre"
$ws.Range("A23").Value = 130886801
$ws.Range("Q23").Value = 434001
$ws.Range("R23").Value = 7052192
$ws.Range("AC23").Value = "Ringhack färska och äldre"
$ws.Range("A24").Value = 130886793
$ws.Range("Q24").Value = 434143
$ws.Range("R24").Value = 7052197
$ws.Range("AC24").Value = "Ringhack äldre"
$ws.Range("A25").Value = 130886794
$ws.Range("Q25").Value = 434140
$ws.Range("R25").Value = 7052192
$ws.Range("AC25").Value = "Ringhack färska och äldre"
$ws.Range("A26").Value = 130886818
$ws.Range("B26").Value = 57884
$ws.Range("E26").Value = 100109
$ws.Range("F26").Value = "Tretåig hackspett"
$ws.Range("G26").Value = "Picoides tridactylus"
$ws.Range("H26").Value = "(Linnaeus, 1758)"
$ws.Range("Q26").Value = 434272
$ws.Range("R26").Value = 7052031
$ws.Range("AC26").Value = "Ringhack äldre"
$ws.Range("A38").Value = 130886835
$ws.Range("B38").Value = 91804
$ws.Range("E38").Value = 1108
$ws.Range("F38").Value = "Harticka"
$ws.Range("G38").Value = "Pelloporus leporinus"
$ws.Range("H38").Value = "(Fr.) Krieglst."
$ws.Range("Q38").Value = 434666
$ws.Range("R38").Value = 7051843
$ws.Range("A39").Value = 130886845
$ws.Range("B39").Value = 91828
$ws.Range("E39").Value = 5432
$ws.Range("F39").Value = "Granticka"
$ws.Range("G39").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H39").Value = ""
$ws.Range("Q39").Value = 434162
$ws.Range("R39").Value = 7052165
$ws.Range("A51").Value = 130886791
$ws.Range("B51").Value = 57884
$ws.Range("E51").Value = 100109
$ws.Range("F51").Value = "Tretåig hackspett"
$ws.Range("G51").Value = "Picoides tridactylus"
$ws.Range("H51").Value = "(Linnaeus, 1758)"
$ws.Range("Q51").Value = 434162
$ws.Range("R51").Value = 7052153
$ws.Range("AC51").Value = "Ringhack"
$ws.Range("A52").Value = 130886765
$ws.Range("Q52").Value = 434505
$ws.Range("R52").Value = 7052005
$ws.Range("A53").Value = 130886788
$ws.Range("Q53").Value = 434171
$ws.Range("R53").Value = 7052213
$ws.Range("AC53").Value = "Ringhack äldre"
$ws.Range("A54").Value = 130886826
$ws.Range("Q54").Value = 434489
$ws.Range("R54").Value = 7051863
$ws.Range("A55").Value = 130886811
$ws.Range("Q55").Value = 434077
$ws.Range("R55").Value = 7052133
$ws.Range("AC55").Value = "Ringhack färska"
$ws.Range("A56").Value = 130886785
$ws.Range("Q56").Value = 434191
$ws.Range("R56").Value = 7052193
$ws.Range("AC56").Value = "Ringhack äldre"
$ws.Range("A57").Value = 130886837
$ws.Range("B57").Value = 91828
$ws.Range("E57").Value = 5432
$ws.Range("F57").Value = "Granticka"
$ws.Range("G57").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H57").Value = ""
$ws.Range("Q57").Value = 434513
$ws.Range("R57").Value = 7052004
$ws.Range("AC57").Value = $null
$ws.Range("A73").Value = 130886832
$ws.Range("B73").Value = 57988
$ws.Range("D73").Value = "LC"
$ws.Range("E73").Value = 103031
$ws.Range("F73").Value = "Lavskrika"
$ws.Range("G73").Value = "Perisoreus infaustus"
$ws.Range("I73").Value = "1"
$ws.Range("K73").Value = ""
$ws.Range("L73").Value = ""
$ws.Range("M73").Value = "födosökande"
$ws.Range("N73").Value = "observerad"
$ws.Range("Q73").Value = 434123
$ws.Range("R73").Value = 7052111
$ws.Range("AC73").Value = $null
$ws.Range("A74").Value = 130886823
$ws.Range("Q74").Value = 434499
$ws.Range("R74").Value = 7051916
$ws.Range("A75").Value = 130886813
$ws.Range("Q75").Value = 434112
$ws.Range("R75").Value = 7052117
$ws.Range("AC75").Value = "Ringhack"
$ws.Range("A76").Value = 130886762
$ws.Range("Q76").Value = 434867
$ws.Range("R76").Value = 7051762
$ws.Range("A77").Value = 130886821
$ws.Range("Q77").Value = 434468
$ws.Range("R77").Value = 7051906
$ws.Range("A78").Value = 130886789
$ws.Range("B78").Value = 57884
$ws.Range("D78").Value = "NT"
$ws.Range("E78").Value = 100109
$ws.Range("F78").Value = "Tretåig hackspett"
$ws.Range("G78").Value = "Picoides tridactylus"
$ws.Range("I78").Value = ""
$ws.Range("K78").Value = $null
$ws.Range("L78").Value = $null
$ws.Range("M78").Value = $null
$ws.Range("N78").Value = $null
$ws.Range("Q78").Value = 434159
$ws.Range("R78").Value = 7052197
$ws.Range("AC78").Value = "Ringhack äldre"

Write-Output "Applied row-permutation edits"